# Refresh cryptocurrency price (D) and 1h volume change (E) figures
# to match the latest coinranking.com snapshot used by the scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.616.02'
$ws.Range('E2').Value = '  -3.87%  '

$ws.Range('D3').Value = '2.399.97'
$ws.Range('E3').Value = '  -4.11%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.19%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '501.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -6.26%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('E8').Value = '  -3.03%  '

$ws.Range('D9').Value = '2.397.13'
$ws.Range('E9').Value = '  -4.28%  '

$ws.Range('E10').Value = '  -3.87%  '

$ws.Range('E11').Value = '  -1.42%  '

$ws.Range('E12').Value = '  -3.62%  '

$ws.Range('E13').Value = '  -10.80%  '

$ws.Range('D14').Value = '2.821.39'
$ws.Range('E14').Value = '  -4.29%  '

$ws.Range('D15').Value = '57.172.20'
$ws.Range('E15').Value = '  -2.65%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.53'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.64%  '

$ws.Range('D18').Value = '2.404.70'
$ws.Range('E18').Value = '  -4.20%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.72%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '309.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.51%  '

$ws.Range('E21').Value = '  -5.45%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.56%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.15%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.50'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.48%  '

$ws.Range('E25').Value = '  +0.30%  '

$ws.Range('D26').Value = '2.491.95'
$ws.Range('E26').Value = '  -4.90%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.373'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.68%  '

$ws.Range('E28').Value = '  -6.17%  '

$ws.Range('E29').Value = '  -3.08%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '174.68'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.59%  '

$ws.Range('E31').Value = '  -4.22%  '

$ws.Range('D32').Value = '0.0₃0710'
$ws.Range('E32').Value = '  -6.05%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.10'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.77%  '

$ws.Range('E34').Value = '  -0.06%  '

$ws.Range('E35').Value = '  -7.36%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.19%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.76'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.87%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.20'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.76%  '

$ws.Range('E39').Value = '  -5.23%  '

$ws.Range('E40').Value = '  -1.89%  '

$ws.Range('E41').Value = '  -6.05%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.769'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.12%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '129.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.30%  '

$ws.Range('E44').Value = '  -4.05%  '

$ws.Range('E45').Value = '  -3.01%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.572'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.97%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '253.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.54%  '

$ws.Range('E48').Value = '  -4.16%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0482'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.46%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.47%  '

$ws.Range('E51').Value = '  -5.35%  '
